$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.740.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.509.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.77"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.508.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("E10").Value = "  +3.94%  "

# Row 11
$ws.Range("E11").Value = "  -0.27%  "

# Row 12
$ws.Range("E12").Value = "  +4.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.32%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.972.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.73%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.690.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.505.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("E19").Value = "  -1.11%  "

# Row 20
$ws.Range("E20").Value = "  -3.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.99%  "

# Row 22
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("E23").Value = "  +1.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.12%  "

# Row 26
$ws.Range("E26").Value = "  -0.35%  "

# Row 27
$ws.Range("E27").Value = "  -0.56%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.647.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0894"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "

# Row 31
$ws.Range("E31").Value = "  +0.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "458.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.12%  "

# Row 33
$ws.Range("E33").Value = "  -2.98%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.88%  "

# Row 37
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "

# Row 39
$ws.Range("E39").Value = "  +0.64%  "

# Row 40
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("E41").Value = "  +1.27%  "

# Row 42
$ws.Range("E42").Value = "  -0.34%  "

# Row 43
$ws.Range("E43").Value = "  +1.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "

# Row 45
$ws.Range("E45").Value = "  -3.45%  "

# Row 46
$ws.Range("E46").Value = "  -6.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.43%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0735"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "

# Row 51
$ws.Range("E51").Value = "  -0.77%  "
